$d = $word.ActiveDocument

$pairs = @(
    @("33×41=", "32×59="),
    @("64×49=", "91×69="),
    @("90×99=", "46×31="),
    @("58×98=", "77×89="),
    @("24×33=", "39×62="),
    @("71×66=", "68×87="),
    @("59×89=", "68×45="),
    @("23×88=", "26×92="),
    @("50×61=", "36×78="),
    @("32×73=", "64×41="),
    @("65×29=", "74×30="),
    @("14×80=", "41×81="),
    @("53×41=", "86×54="),
    @("15×28=", "66×93="),
    @("84×43=", "47×65="),
    @("34×56=", "11×68="),
    @("80×61=", "65×75="),
    @("77×50=", "57×25="),
    @("38×14=", "31×45="),
    @("46×45=", "45×74="),
    @("53×46=", "53×64="),
    @("59×36=", "62×51="),
    @("94×69=", "42×95="),
    @("76×27=", "23×53="),
    @("18×62=", "50×26=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
